$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 115 <= old Row 116 (B:AD)
$ws.Range("B115").Value = 7919322
$ws.Range("C115").Value = "Estonia Meistriliiga"
$ws.Range("D115").Value = 45360.39583333334
$ws.Range("E115").Value = "FC Kuressaare"
$ws.Range("F115").Value = "FC Levadia Tallinn"
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 6
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 2
$ws.Range("K115").Value = "A"
$ws.Range("L115").Value = 11
$ws.Range("M115").Value = 6
$ws.Range("N115").Value = 1.166
$ws.Range("O115").Value = 15
$ws.Range("P115").Value = 8.5
$ws.Range("Q115").Value = 1.125
$ws.Range("R115").Value = 2.5
$ws.Range("S115").Value = 1.825
$ws.Range("T115").Value = 1.975
$ws.Range("U115").Value = 3.25
$ws.Range("V115").Value = 1.9
$ws.Range("W115").Value = 1.9
$ws.Range("X115").Value = -1
$ws.Range("Y115").Value = -1
$ws.Range("Z115").Value = 0.125
$ws.Range("AA115").Value = -1
$ws.Range("AB115").Value = 0.9750000000000001
$ws.Range("AC115").Value = 0.8999999999999999
$ws.Range("AD115").Value = -1

# Row 116 <= old Row 115 (B:AD)
$ws.Range("B116").Value = 7919323
$ws.Range("C116").Value = "Estonia Meistriliiga"
$ws.Range("D116").Value = 45360.39583333334
$ws.Range("E116").Value = "JK Nomme Kalju"
$ws.Range("F116").Value = "JK Trans Narva"
$ws.Range("G116").Value = 3
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 2
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = "H"
$ws.Range("L116").Value = 1.285
$ws.Range("M116").Value = 5.5
$ws.Range("N116").Value = 6.5
$ws.Range("O116").Value = 1.571
$ws.Range("P116").Value = 4.75
$ws.Range("Q116").Value = 4.2
$ws.Range("R116").Value = -1
$ws.Range("S116").Value = 1.925
$ws.Range("T116").Value = 1.875
$ws.Range("U116").Value = 2.75
$ws.Range("V116").Value = 1.875
$ws.Range("W116").Value = 1.925
$ws.Range("X116").Value = 0.571
$ws.Range("Y116").Value = -1
$ws.Range("Z116").Value = -1
$ws.Range("AA116").Value = 0.925
$ws.Range("AB116").Value = -1
$ws.Range("AC116").Value = 0.4375
$ws.Range("AD116").Value = -0.5

# Row 120 <= old Row 121 (B:AD)
$ws.Range("B120").Value = 7721087
$ws.Range("C120").Value = "Estonia Meistriliiga"
$ws.Range("D120").Value = 45367.39583333334
$ws.Range("E120").Value = "Paide Linnameeskond"
$ws.Range("F120").Value = "FC Flora Tallinn"
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = 1
$ws.Range("I120").Value = 1
$ws.Range("J120").Value = 1
$ws.Range("K120").Value = "H"
$ws.Range("L120").Value = 2.2
$ws.Range("M120").Value = 3.3
$ws.Range("N120").Value = 2.8
$ws.Range("O120").Value = 1.85
$ws.Range("P120").Value = 3.6
$ws.Range("Q120").Value = 3.4
$ws.Range("R120").Value = -0.5
$ws.Range("S120").Value = 1.9
$ws.Range("T120").Value = 1.9
$ws.Range("U120").Value = 2.5
$ws.Range("V120").Value = 1.95
$ws.Range("W120").Value = 1.85
$ws.Range("X120").Value = 0.8500000000000001
$ws.Range("Y120").Value = -1
$ws.Range("Z120").Value = -1
$ws.Range("AA120").Value = 0.8999999999999999
$ws.Range("AB120").Value = -1
$ws.Range("AC120").Value = 0.95
$ws.Range("AD120").Value = -1

# Row 121 <= old Row 120 (B:AD)
$ws.Range("B121").Value = 7721007
$ws.Range("C121").Value = "Estonia Meistriliiga"
$ws.Range("D121").Value = 45367.39583333334
$ws.Range("E121").Value = "JK Trans Narva"
$ws.Range("F121").Value = "JK Tammeka Tartu"
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 5
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2
$ws.Range("K121").Value = "A"
$ws.Range("L121").Value = 2.25
$ws.Range("M121").Value = 3.3
$ws.Range("N121").Value = 2.75
$ws.Range("O121").Value = 2.1
$ws.Range("P121").Value = 3.25
$ws.Range("Q121").Value = 3
$ws.Range("R121").Value = -0.25
$ws.Range("S121").Value = 1.875
$ws.Range("T121").Value = 1.925
$ws.Range("U121").Value = 2.5
$ws.Range("V121").Value = 1.825
$ws.Range("W121").Value = 1.975
$ws.Range("X121").Value = -1
$ws.Range("Y121").Value = -1
$ws.Range("Z121").Value = 2
$ws.Range("AA121").Value = -1
$ws.Range("AB121").Value = 0.925
$ws.Range("AC121").Value = 0.825
$ws.Range("AD121").Value = -1

# Row 175 odds updates
$ws.Range("O175").Value = 4.75
$ws.Range("P175").Value = 4.1
$ws.Range("Q175").Value = 1.55
$ws.Range("R175").Value = 1
$ws.Range("S175").Value = 1.825
$ws.Range("T175").Value = 1.975
